# "Draft of generalized LCOE with future costs" (#373)
#
# Adds a new "FFY" (Final Future Year) sheet after "FY", mirroring the
# existing "FY" sheet's year-subscript pattern but anchored on a new
# "Final Time" input on the "About" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "FFY" worksheet as the last tab -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ffy = $wb.Worksheets.Add($null, $lastSheet)
$ffy.Name = "FFY"

# Mirrors FY!A1 ("FutureYear") / FY!A2 ("Year"&About!A8) but for the new
# "Final Time" input cell (About!A9).
$ffy.Range("A1").Value = "Final Future Year"
$ffy.Range("A2").Formula = "=""Year""&About!A9"

# --- 2. Add the "Final Time" input row to the "About" sheet ---------------
$about = $wb.Worksheets.Item("About")

$about.Range("A9").Value = 2050
$about.Range("A9").Font.Bold = $true
$about.Range("A9").Interior.Color = 65535

$about.Range("B9").Value = "Final Time"
$about.Range("B9").Font.Bold = $true

# --- 3. Restore per-sheet selections / active sheet ------------------------
$fy = $wb.Worksheets.Item("FY")
$fy.Range("A3").Select()

$about.Range("A10").Select()

# Selecting FFY last makes it the active/visible tab, matching the diff.
$ffy.Range("A2").Select()
